$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 3 and 4 (gregneri12/halachme and armonravid2/armonravid
# reviews). The old row 5 (veredsnir12/kevinkors122) shifts up and becomes
# the new row 3.
$ws.Rows("3:4").Delete() | Out-Null

# Hyperlinks.Add() always reformats its target cell with the built-in
# "Hyperlink" style, so stash the current (correct) formatting of the four
# cells that will get new hyperlinks in a scratch area first.
$ws.Range("C2").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy()
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy()
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy()
$ws.Range("H4").PasteSpecial(-4122) | Out-Null

# Rebuild the hyperlinks collection so it only references the remaining
# two data rows (the collection doesn't support removing single items).
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:snizzvered@gmail.com", "", "", "snizzvered@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:krigelron@gmail.com", "", "", "krigelron@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com") | Out-Null

# Restore the original cell formatting that Hyperlinks.Add() overwrote.
$ws.Range("H1").Copy()
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy()
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").Copy()
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Copy()
$ws.Range("D3").PasteSpecial(-4122) | Out-Null

# Clean up the scratch cells.
$ws.Range("H1:H4").Clear() | Out-Null

# Move the active selection to A3, matching the saved view state.
$ws.Range("A3").Select() | Out-Null
